$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data_smooth")

$rng = $ws.Range("K4:O4")

$ws.Range("K4").Value = 400.18799999999999
$ws.Range("L4").Value = 413.70299999999997
$ws.Range("M4").Value = 399.68799999999999
$ws.Range("N4").Value = 385.375
$ws.Range("O4").Value = 395.03100000000001

# The pasted-in values came without the thin top border the row used to
# carry (fill stays the existing yellow), so drop the border on the range.
$rng.Borders.LineStyle = -4142

$rng.Select()
